$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shared-string text for row 31 (day 30 draft entry)
$ayats = 'Surah Al Imran, 15 - 25'
$tags = 'Five commandments, People beloved to Allah, Duty of a muslim, Response to animosity'
$content = @'
h1: Purpose of these blogs?
p.note: I will try to stick to Quran most part of this blog. 
quote: Say, "Shall I inform you of [something] better than that? For those who fear Allah will be gardens in the presence of their Lord beneath which rivers flow, wherein they abide eternally, and purified spouses and approval from Allah . And Allah is Seeing of [His] servants. Those who say, "Our Lord, indeed we have believed, so forgive us our sins and protect us from the punishment of the Fire”. The patient, the true, the obedient, those who spend [in the way of Allah ], and those who seek forgiveness before dawn. <br> - Surah Al Imran verse 15 - 17
p: A beautiful life is coming ahead. Gardens are green, rivers are flowing, calmness dominates, happiness lasts eternally, with spouses loving us and above all the creator providing for us all visible in its entirety. If we doubt these gardens, we really are not enjoying our life. We are doubting our life, present and future. Allah is watching us and he is protecting only those among us who repent to sins and ask for His protection from the fire. 
p: Right after repentance, we are immediately expected to:-
p.b-left: <b>1. Be patient.</b> Not impulsive and reactive but disciplined and maturely.
p.b-left: <b>2. Be truthful,</b> to our spouses, friends, colleagues and subordinates.
p.b-left: <b>3. Be obedient.</b> We do not resort to intoxicants to find comfort. We avoid music, lies and everything wrong to follow His commandment. 
p.b-left: <b>4. Spend in His way.</b> Parents first, then relatives, orphans, needy and travellers; as mentioned in following verse:-
quote: They ask you, [O Muhammad], what they should spend. Say, "Whatever you spend of good is [to be] for parents and relatives and orphans and the needy and the traveler. And whatever you do of good - indeed, Allah is Knowing of it." <br> - Surah Baqarah verse 215
p.b-left: <b>5. Repent before dawn (Tahajjud).</b> Early in the morning when we believe our Creator is closest to us waiting for us to ask Him, that He gives at once.
p: Purpose of these blogs is to build these traits in myself. I have to work hard to secure my death’s comfort. 
h3: Next purpose.. Write non stop with all my love.
quote: Indeed, those who have believed and those who have emigrated and fought in the cause of Allah - those expect the mercy of Allah . And Allah is Forgiving and Merciful. <br> - Surah Al Imran verse 218
p: After following traits and wanting forgiveness of Allah (swt) I have to strive harder, farther and longer. Longing to expect the mercy of Allah (swt). If I am worried He will not forgive me, I am wrong. He is forgiving and merciful, He will forgive me. He loves me today. Today He will bless me with comfort. I have to trust Him and write whole heartedly, lovely and consistently. Every night I have to defeat my nafs and fight in the cause of Allah. Every night He tells me, He is forgiving and merciful.
h3: Do not act macho man.
quote: They ask you about wine and gambling. Say, "In them is great sin and [yet, some] benefit for people. But their sin is greater than their benefit." And they ask you what they should spend. Say, "The excess [beyond needs]." Thus Allah makes clear to you the verses [of revelation] that you might give thought. <br> - Surah Al Imran verse 219
p: If people say, I am having trouble enjoying my life. If they say I am not drinking wine and gambling more, I am not really enjoying my life. I can not give an ear to what they say. If someone says “I can not really follow Prophet (PBUH), considering him being a Super Human. I can not really spend my hard earned money on his teachings.”. I just have to write the response here, “You should spend the excess.” 
p: I do not go all macho man justifying my comments and words. I just have to write what Quran says. 
h3: Reward is near..
quote: Allah witnesses that there is no deity except Him, and [so do] the angels and those of knowledge - [that He is] maintaining [creation] in justice. There is no deity except Him, the Exalted in Might, the Wise. <br> - Surah Al Imran verse 18
p: Allah is maintainer of justice. Justice that is visible to His angels and His creation with knowledge and understanding. Justly sooner than ever we will be rewarded for our efforts. 
h3: Quickly embrace Quran
quote: Indeed, the religion in the sight of Allah is Islam. And those who were given the Scripture did not differ except after knowledge had come to them - out of jealous animosity between themselves. And whoever disbelieves in the verses of Allah, then indeed, Allah is swift in [taking] account. <br> - Surah Al Imran verse 19
p: I can not be worried about what to write. What really will impact the reader. People who were given the scripture differed out of animosity. They really understood the Quran and the Prophet (PBUH), but they failed to embrace it. 
h3: Coming to the point..
quote: So if they argue with you, say, "I have submitted myself to Allah [in Islam], and [so have] those who follow me." And say to those who were given the Scripture and [to] the unlearned, "Have you submitted yourselves?" And if they submit [in Islam], they are rightly guided; but if they turn away - then upon you is only the [duty of] notification. And Allah is Seeing of [His] servants. <br> - Surah Al Imran verse 20
p: The purpose of these blogs is mentioned in above ayat. It is to read Quran, act on it, share it and avoid arguments. Purpose is to notify / deliver Quran on internet. May be there are excuses like Islamic websites are not really designed well. There are some but they do not really talk much. May be these blogs fill this gap and notify people. 
p: Yet if these blogs did not create the impact, I have to keep on walking patiently, truthfuly, obediently on the way of Allah (swt). I have to keep walking, head lowered towards this light that is going to meet me on the next junction. These blogs are actually Quran talking to us. Nothing more than this, nothing less. Plain old Quran talking to a group of people..
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
'@

# Populate row 31 (Ser 30 / 30-Jan-2020 draft)
# Shared-string table order matters: Ayats, then Tags, then the long Content
# blog post, matching the order new <si> entries were appended upstream.
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 43860
$ws.Range("B31").NumberFormat = $ws.Range("B30").NumberFormat
$ws.Range("C31").Value = $ayats
$ws.Range("F31").Value = $tags
$ws.Range("D31").Value = $content
$ws.Range("E31").Value = "Qasim Ali"

# Row grows to the max row height to fit the long blog entry, like the rows above it
$ws.Rows.Item(31).RowHeight = 409.6

# Move the saved cursor/selection down onto the newly-filled row
$ws.Activate() | Out-Null
$ws.Range("D31").Select() | Out-Null
